# Remove the three "BW" (back-wait?) captain rows that were merged/removed,
# and add Tsawwassen Tuesday context per the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete whole rows (from bottom to top so row numbers of earlier rows stay valid)
# Row 38: University BW / Lori Hutchinson
# Row 36: Tsawwassen Springs BW / Moira Milligan
# Row 30: Quilchena BW / Casandra Kobayashi
$ws.Rows.Item(38).Delete()
$ws.Rows.Item(36).Delete()
$ws.Rows.Item(30).Delete()

$ws.Range("A30").Select()
